# Restore C10 ("From" value of rule R30 on the Rules sheet) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
